$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.476.48"
$ws.Range("E2").Value = "  +5.04%  "

# Row 3
$ws.Range("D3").Value = "1.724.43"
$ws.Range("E3").Value = "  +4.13%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.88"
$ws.Range("E5").Value = "  +3.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5377"
$ws.Range("E6").Value = "  +2.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2683"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06597"
$ws.Range("E9").Value = "  +3.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.66"
$ws.Range("E10").Value = "  +5.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.636"
$ws.Range("E12").Value = "  +0.37%  "

# Row 13
$ws.Range("D13").Value = "1.722.60"
$ws.Range("E13").Value = "  +3.95%  "

# Row 14
$ws.Range("D14").Value = "1.958.81"
$ws.Range("E14").Value = "  +3.97%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5871"
$ws.Range("E15").Value = "  +4.60%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8275"
$ws.Range("E16").Value = "  +1.14%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.00"
$ws.Range("E17").Value = "  +3.81%  "

# Row 18
$ws.Range("D18").Value = "27.496.03"
$ws.Range("E18").Value = "  +5.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.54"
$ws.Range("E19").Value = "  +14.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.09%  "

# Row 21
$ws.Range("E21").Value = "  +1.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("E22").Value = "  +1.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.095"
$ws.Range("E23").Value = "  +2.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.62"
$ws.Range("E25").Value = "  +1.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1232"
$ws.Range("E26").Value = "  +3.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.689"
$ws.Range("E27").Value = "  +11.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.409"
$ws.Range("E28").Value = "  +2.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.68"
$ws.Range("E29").Value = "  +4.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05553"
$ws.Range("E30").Value = "  +1.34%  "

# Row 31
$ws.Range("E31").Value = "  +2.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.534"
$ws.Range("E32").Value = "  +2.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.462"
$ws.Range("E33").Value = "  +2.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.656"
$ws.Range("E34").Value = "  +6.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.450"
$ws.Range("E35").Value = "  +1.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9591"
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.813"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5919"
$ws.Range("E38").Value = "  +4.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01640"
$ws.Range("E39").Value = "  +3.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.862"
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8555"
$ws.Range("E41").Value = "  +2.91%  "

# Row 42
$ws.Range("D42").Value = "1.053.81"
$ws.Range("E42").Value = "  +2.54%  "

# Row 43
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.66"
$ws.Range("E44").Value = "  +0.74%  "

# Row 45
$ws.Range("D45").Value = "1.866.51"
$ws.Range("E45").Value = "  +4.00%  "

# Row 46
$ws.Range("E46").Value = "  +8.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.88"
$ws.Range("E47").Value = "  +1.65%  "

# Row 48 (swap: now EnergySwap data)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.195"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49 (swap: now Mantle data)
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4441"
$ws.Range("E49").Value = "  +2.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9993"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05276"
$ws.Range("E51").Value = "  +1.28%  "

